$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newSamplesQuery = @'
SELECT
    smp.sample_id AS "Sample ID",
    prt.participant_id AS "Participant ID",
    std.dbgap_accession AS "Study ID",
    smp.anatomic_site AS "Sample Anatomic Site",
    COALESCE(CASE WHEN smp.participant_age_at_collection = -999 THEN 'Not Reported' ELSE smp.participant_age_at_collection END, 0) AS "Age at Sample Collection (days)",
    COALESCE(smp.sample_tumor_status, '') AS "Sample Tumor Status",
    COALESCE(smp.tumor_classification, '') AS "Sample Tumor Classification",
    Null  AS "Sample Diagnosis"
FROM 
    df_study std
LEFT JOIN 
    df_participant prt ON std.id = prt."study.id"
LEFT JOIN 
    df_sample smp ON prt.id = smp."participant.id"
LEFT JOIN 
    df_diagnosis dgn ON prt.id = dgn."participant.id"
WHERE 
    std.dbgap_accession = 'phs002371' 
    AND prt.sex_at_birth = 'Male'
	and dgn.anatomic_site= 'C42.0 : Blood'
ORDER BY 
    smp.sample_id ASC;
'@

# Row 4 (SamplesTab) now holds the expanded query with extra sample columns.
# Row 5 (FilesTab) keeps its existing query text (file_data query) unchanged.
$ws.Cells.Item(4, 2).Value = $newSamplesQuery
$ws.Rows.Item(4).RowHeight = 393.75
$ws.Rows.Item(5).RowHeight = 409.5

# Update the active view/selection to match the saved state
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("B4").Select()
